$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "admin#admin"
$ws.Range("A1").Value = "Admin"
